# Apply the "Contact Tracing Data" edits described by the diff:
#  - Update row 2 values (name casing, address, time, date, vaccine, symptoms)
#  - Delete row 3 entirely (it was a duplicate entry)
#  - Widen columns F and H
#  - Reset the sheet view (remove frozen/top-left scroll, change selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 cell values ---
$ws.Range("A2").Value = "Mary Grace"
$ws.Range("B2").Value = "Dellomos"
$ws.Range("E2").Value = "etivac"
$ws.Range("F2").Value = "4:50 A.M."
$ws.Range("G2").Value = "26/07/2023"
$ws.Range("H2").Value = "First Booster Shot"
$ws.Range("I2").Value = "No, Fever, Cough"

# --- Delete row 3 (duplicate entry), shifting dimension back to A1:I2 ---
$ws.Rows.Item(3).Delete()

# --- Column widths ---
# (the host quantizes ColumnWidth to 1/6-character steps, so the input is
#  chosen to land as close as possible to the target stored widths of
#  12.6640625 and 21.109375 characters)
$ws.Columns.Item(6).ColumnWidth = 11.833333333333332
$ws.Columns.Item(8).ColumnWidth = 20.333333333333336

# --- Sheet view / selection updates ---
$ws.Activate()
$ws.Range("E5").Select()
